# Applies the "output generated at 456a3b4" data refresh to 上海-漫展信息.xlsx
# (updates "想去人数" / "最低票价" counters across all four sheets)
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F2").Value = 157
$ws.Range("F3").Value = 174
$ws.Range("F4").Value = 2205
$ws.Range("F5").Value = 4351
$ws.Range("F6").Value = 580
$ws.Range("F8").Value = 1350
$ws.Range("F9").Value = 670
$ws.Range("F10").Value = 384
$ws.Range("F11").Value = 108
$ws.Range("F13").Value = 677543
$ws.Range("F15").Value = 581
$ws.Range("F16").Value = 1492
$ws.Range("F17").Value = 676
$ws.Range("F19").Value = 1308
$ws.Range("F20").Value = 2315
$ws.Range("F21").Value = 1164
$ws.Range("F22").Value = 2729
$ws.Range("F24").Value = 890
$ws.Range("F25").Value = 1582
$ws.Range("F28").Value = 1023
$ws.Range("F29").Value = 1106
$ws.Range("F31").Value = 94
$ws.Range("F32").Value = 2057
$ws.Range("F34").Value = 1344
$ws.Range("F35").Value = 3091
$ws.Range("F37").Value = 1155
$ws.Range("F38").Value = 50
$ws.Range("F40").Value = 2635
$ws.Range("F43").Value = 3183
$ws.Range("F47").Value = 165
$ws.Range("F48").Value = 674
$ws.Range("F49").Value = 26

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")

$ws.Range("G5").Value = "不可售"
$ws.Range("F9").Value = 115
$ws.Range("F11").Value = 145219
$ws.Range("G11").Value = 380
$ws.Range("F12").Value = 145219
$ws.Range("G12").Value = 380
$ws.Range("F13").Value = 12
$ws.Range("F17").Value = 98
$ws.Range("F18").Value = 235
$ws.Range("F19").Value = 343
$ws.Range("F21").Value = 425
$ws.Range("F22").Value = 181
$ws.Range("F23").Value = 87
$ws.Range("F24").Value = 97
$ws.Range("F26").Value = 655
$ws.Range("F29").Value = 13
$ws.Range("F31").Value = 375
$ws.Range("F32").Value = 282
$ws.Range("F39").Value = 113
$ws.Range("F41").Value = 190
$ws.Range("F43").Value = 10

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")

$ws.Range("F4").Value = 3140
$ws.Range("F5").Value = 255
$ws.Range("F7").Value = 838
$ws.Range("F8").Value = 1229
$ws.Range("F9").Value = 650
$ws.Range("F10").Value = 1618
$ws.Range("F11").Value = 146
$ws.Range("F12").Value = 2040

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F2").Value = 838
$ws.Range("F3").Value = 650
$ws.Range("F4").Value = 157
$ws.Range("F5").Value = 1618
$ws.Range("F6").Value = 174
$ws.Range("F7").Value = 2205
$ws.Range("F8").Value = 146
$ws.Range("F9").Value = 2040
$ws.Range("F10").Value = 4352
$ws.Range("F11").Value = 580
$ws.Range("F12").Value = 1350
$ws.Range("F13").Value = 670
$ws.Range("F14").Value = 384
$ws.Range("F16").Value = 677552
$ws.Range("F17").Value = 115
$ws.Range("F20").Value = 145220
$ws.Range("G20").Value = 380
$ws.Range("F21").Value = 1492
$ws.Range("F22").Value = 676
$ws.Range("F24").Value = 1308
$ws.Range("F25").Value = 2316
$ws.Range("F26").Value = 1164
$ws.Range("F27").Value = 2729
$ws.Range("F29").Value = 890
$ws.Range("F31").Value = 1582
$ws.Range("F33").Value = 181
$ws.Range("F35").Value = 1025
$ws.Range("F36").Value = 1106
$ws.Range("F37").Value = 94
$ws.Range("F38").Value = 2057
$ws.Range("F39").Value = 1344
$ws.Range("F40").Value = 3091
$ws.Range("F42").Value = 1155
$ws.Range("F43").Value = 375
$ws.Range("F44").Value = 282
$ws.Range("F46").Value = 2635
$ws.Range("F49").Value = 3183
$ws.Range("F52").Value = 165
$ws.Range("F53").Value = 674
$ws.Range("F54").Value = 26
